$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily row appended to the log (row 71).
# Force column A to be stored as literal text ("2025/10/07") rather than
# being auto-converted to a date serial number, matching the original
# sheet's convention (all date cells are inline/shared strings, not real
# dates). Temporarily apply a text number format while assigning the
# value, then restore the default "Normal" style so the new row's cells
# end up with the same (unstyled) formatting as the other data rows.
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2025/10/07"
$ws.Range("A71").Style = "Normal"

$ws.Range("B71").Value = "火"
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = 5
